# Apply updated dSF (column F) values per the commit:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -4
$ws.Range("F9").Value = -6
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = -2
$ws.Range("F20").Value = 0
$ws.Range("F22").Value = 9
$ws.Range("F27").Value = -1
$ws.Range("F28").Value = -2
$ws.Range("F33").Value = -7
$ws.Range("F35").Value = 3
